$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the two "Annual Population Survey" link cells to note they are
# direct downloads (append " (direct download)" before the closing </a> tag).
$ws.Range("B2").Value = "<a href='https://www.nomisweb.co.uk/livelinks/16244.xlsx'>Annual Population Survey (direct download)</a>"
$ws.Range("B4").Value = "<a href='https://www.nomisweb.co.uk/livelinks/16243.xlsx'>Annual Population Survey (direct download)</a>"

# Reflect the final selection state left in the saved file.
$ws.Range("B14").Select()
